# "break out stock.yaml completed"
# - Append 9 new rows (435-443) of stock data to the "day" sheet.
# - On the "week" sheet, rows 158-169 column D (bsecode) values are
#   rewritten as numeric cells instead of text cells (same displayed value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "day" sheet: append new rows 435-443
# ---------------------------------------------------------------------------
$day = $wb.Worksheets.Item("day")

$newRows = @(
    @(1, "PAGEIND",    "Page Industries Limited",    532827, 1.17,               42000,   16114,   "26/08/2024 11:35:36"),
    @(2, "SHREECEM",   "Shree Cements Limited",      500387, 0.79,               24900,   17218,   "26/08/2024 11:35:36"),
    @(3, "ABB",        "Abb India Limited",          500002, -0.07000000000000001, 7790.4, 189412,  "26/08/2024 11:35:36"),
    @(4, "PIIND",      "Pi Industries Limited",      523642, -0.46,              4399.95, 344105,  "26/08/2024 11:35:36"),
    @(5, "VOLTAS",     "Voltas Limited",             500575, 1.8,                1721.05, 1983382, "26/08/2024 11:35:36"),
    @(6, "UNITDSPR",   "United Spirits Ltd",         532432, 1.94,               1457.45, 845392,  "26/08/2024 11:35:36"),
    @(7, "ZYDUSLIFE",  "Zydus Lifesciences Ltd",     532321, -5.96,              1108.45, 7365014, "26/08/2024 11:35:36"),
    @(8, "SUNTV",      "Sun Tv Network Limited",     532733, 2.96,               792.15,  1126572, "26/08/2024 11:35:36"),
    @(9, "BERGEPAINT", "Berger Paints (i) Limited",  509480, 1.18,               572.35,  810296,  "26/08/2024 11:35:36")
)

$startRow = 435
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $day.Cells.Item($r, 1).Value = $row[0]
    $day.Cells.Item($r, 2).Value = $row[1]
    $day.Cells.Item($r, 3).Value = $row[2]
    $day.Cells.Item($r, 4).Value = $row[3]
    $day.Cells.Item($r, 5).Value = $row[4]
    $day.Cells.Item($r, 6).Value = $row[5]
    $day.Cells.Item($r, 7).Value = $row[6]
    $day.Cells.Item($r, 8).Value = "day"
    $day.Cells.Item($r, 9).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. "week" sheet: rows 158-169, column D (bsecode) becomes numeric
# ---------------------------------------------------------------------------
$week = $wb.Worksheets.Item("week")

$bsecodes = @{
    158 = 500387
    159 = 540005
    160 = 500820
    161 = 533150
    162 = 533309
    163 = 533273
    164 = 500112
    165 = 500253
    166 = 500850
    167 = 532400
    168 = 507685
    169 = 531213
}

foreach ($r in $bsecodes.Keys) {
    $week.Cells.Item($r, 4).Value = $bsecodes[$r]
}

Write-Output "edit applied"
